# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 16 becomes the new worker (ANDRESON EXMIT ZUÑIGA ARISTIZABAL, doc 1201256725,
# periodo 1802, valor mora 72000, salario basico 1800000) - this was previously
# at row 23.
# Rows 17-23 now all belong to KEVIN BARRIOS OSORIO (doc 1143373237), with the
# "Periodo Mora" column running sequentially 2205..2211 (row 23 keeps the
# 37800 / 1350000 pair that used to sit on row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: ANDRESON EXMIT ZUÑIGA ARISTIZABAL
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1201256725"
$ws.Range("D16").Value = "ANDRESON EXMIT ZUÑIGA ARISTIZABAL"
$ws.Range("E16").Value = "1802"
$ws.Range("F16").Value = 72000
$ws.Range("G16").Value = 1800000

# Rows 17-22: KEVIN BARRIOS OSORIO, periodos 2205-2210, valor 54000 / salario 1350000
$periodos = @("2205", "2206", "2207", "2208", "2209", "2210")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 17 + $i
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1143373237"
    $ws.Range("D$row").Value = "KEVIN BARRIOS OSORIO"
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = 54000
    $ws.Range("G$row").Value = 1350000
}

# Row 23: KEVIN BARRIOS OSORIO, periodo 2211, valor 37800 / salario 1350000
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143373237"
$ws.Range("D23").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E23").Value = "2211"
$ws.Range("F23").Value = 37800
$ws.Range("G23").Value = 1350000
